# Apply the "Updated cryptos list" data refresh to the Price (D) and
# Volume(1h) (E) columns of the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values look like plain numbers (e.g. "1.000", "0.9999").
# Excel would normally reinterpret/round such text when assigned through
# the COM object model, so those specific cells are pre-formatted as text
# before the value is written, preserving the exact displayed digits.
$priceTextCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D19", "D21", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D46", "D48", "D50", "D51")
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Price (column D) updates
$ws.Range("D2").Value = "27.149.05"
$ws.Range("D3").Value = "1.867.75"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D5").Value = "306.19"
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").Value = "0.5159"
$ws.Range("D9").Value = "0.07149"
$ws.Range("D10").Value = "0.8920"
$ws.Range("D11").Value = "20.75"
$ws.Range("D12").Value = "0.07547"
$ws.Range("D13").Value = "1.863.57"
$ws.Range("D14").Value = "5.311"
$ws.Range("D15").Value = "89.57"
$ws.Range("D16").Value = "1.000"
$ws.Range("D17").Value = "0.000008485"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D20").Value = "27.188.09"
$ws.Range("D21").Value = "4.997"
$ws.Range("D22").Value = "2.100.32"
$ws.Range("D24").Value = "6.456"
$ws.Range("D25").Value = "1.837"
$ws.Range("D26").Value = "145.59"
$ws.Range("D27").Value = "17.97"
$ws.Range("D28").Value = "2.086"
$ws.Range("D29").Value = "112.91"
$ws.Range("D30").Value = "4.661"
$ws.Range("D31").Value = "4.683"
$ws.Range("D32").Value = "0.09252"
$ws.Range("D33").Value = "0.05123"
$ws.Range("D35").Value = "1.161"
$ws.Range("D36").Value = "0.7252"
$ws.Range("D37").Value = "0.02034"
$ws.Range("D38").Value = "3.105"
$ws.Range("D39").Value = "2.505"
$ws.Range("D41").Value = "0.5289"
$ws.Range("D42").Value = "6.512"
$ws.Range("D43").Value = "116.55"
$ws.Range("D44").Value = "8.318"
$ws.Range("D46").Value = "0.9998"
$ws.Range("D48").Value = "9.987"
$ws.Range("D50").Value = "36.65"
$ws.Range("D51").Value = "63.63"

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("E21").Value = "  -2.71%  "
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("E23").Value = "  -3.37%  "
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("E26").Value = "  -5.16%  "
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("E28").Value = "  -3.32%  "
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("E32").Value = "  +2.55%  "
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("E34").Value = "  -3.24%  "
$ws.Range("E35").Value = "  -5.61%  "
$ws.Range("E36").Value = "  -7.04%  "
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("E41").Value = "  -4.87%  "
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("E47").Value = "  -3.92%  "
$ws.Range("E48").Value = "  -4.82%  "
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("E51").Value = "  -4.66%  "
